# PowerShell Excel COM-interop script
# Commit: "primera version LSA Y TFIDF"
#
# Changes applied:
#  1. Delete worksheet "Cluster4" (sheetId 5) -- cluster set shrinks from 5 to 4 clusters.
#  2. Topic columns are renumbered/extended: Topic_0..Topic_15 (16 topics, cols C:R)
#     becomes Topic_1..Topic_19 (19 topics, cols C:U) on every remaining sheet.
#  3. Row data (cluster membership / topic-weight distributions) is refreshed with
#     the new TF-IDF/LSA values, including a couple of rows moving between sheets,
#     two brand new rows appended to "Cluster2", and one row removed from "Cluster0".

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# --- Remove the "Cluster4" worksheet entirely ---
$wb.Worksheets.Item("Cluster4").Delete() | Out-Null

# --- Update worksheet "Cluster0" ---
$ws = $wb.Worksheets.Item("Cluster0")
$ws.Rows.Item(4).Delete() | Out-Null
$ws.Range("B1").Value = "clusters"
$ws.Range("C1").Value = "Topic_1"
$ws.Range("D1").Value = "Topic_2"
$ws.Range("E1").Value = "Topic_3"
$ws.Range("F1").Value = "Topic_4"
$ws.Range("G1").Value = "Topic_5"
$ws.Range("H1").Value = "Topic_6"
$ws.Range("I1").Value = "Topic_7"
$ws.Range("J1").Value = "Topic_8"
$ws.Range("K1").Value = "Topic_9"
$ws.Range("L1").Value = "Topic_10"
$ws.Range("M1").Value = "Topic_11"
$ws.Range("N1").Value = "Topic_12"
$ws.Range("O1").Value = "Topic_13"
$ws.Range("P1").Value = "Topic_14"
$ws.Range("Q1").Value = "Topic_15"
$ws.Range("R1").Value = "Topic_16"
$ws.Range("S1").Value = "Topic_17"
$ws.Range("T1").Value = "Topic_18"
$ws.Range("U1").Value = "Topic_19"
$ws.Range("A2").Value = "antena3_2019 09 17_morning_new"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.428
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.5649999999999999
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("A3").Value = "antena3_2019 09 17_afternoon_new"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0

# --- Update worksheet "Cluster1" ---
$ws = $wb.Worksheets.Item("Cluster1")
$ws.Range("B1").Value = "clusters"
$ws.Range("C1").Value = "Topic_1"
$ws.Range("D1").Value = "Topic_2"
$ws.Range("E1").Value = "Topic_3"
$ws.Range("F1").Value = "Topic_4"
$ws.Range("G1").Value = "Topic_5"
$ws.Range("H1").Value = "Topic_6"
$ws.Range("I1").Value = "Topic_7"
$ws.Range("J1").Value = "Topic_8"
$ws.Range("K1").Value = "Topic_9"
$ws.Range("L1").Value = "Topic_10"
$ws.Range("M1").Value = "Topic_11"
$ws.Range("N1").Value = "Topic_12"
$ws.Range("O1").Value = "Topic_13"
$ws.Range("P1").Value = "Topic_14"
$ws.Range("Q1").Value = "Topic_15"
$ws.Range("R1").Value = "Topic_16"
$ws.Range("S1").Value = "Topic_17"
$ws.Range("T1").Value = "Topic_18"
$ws.Range("U1").Value = "Topic_19"
$ws.Range("A2").Value = "antena3_2019 09 16_morning_new"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0.988
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("A3").Value = "antena3_2019 09 18_morning_new"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0

# --- Update worksheet "Cluster2" ---
$ws = $wb.Worksheets.Item("Cluster2")
$ws.Range("B1").Value = "clusters"
$ws.Range("C1").Value = "Topic_1"
$ws.Range("D1").Value = "Topic_2"
$ws.Range("E1").Value = "Topic_3"
$ws.Range("F1").Value = "Topic_4"
$ws.Range("G1").Value = "Topic_5"
$ws.Range("H1").Value = "Topic_6"
$ws.Range("I1").Value = "Topic_7"
$ws.Range("J1").Value = "Topic_8"
$ws.Range("K1").Value = "Topic_9"
$ws.Range("L1").Value = "Topic_10"
$ws.Range("M1").Value = "Topic_11"
$ws.Range("N1").Value = "Topic_12"
$ws.Range("O1").Value = "Topic_13"
$ws.Range("P1").Value = "Topic_14"
$ws.Range("Q1").Value = "Topic_15"
$ws.Range("R1").Value = "Topic_16"
$ws.Range("S1").Value = "Topic_17"
$ws.Range("T1").Value = "Topic_18"
$ws.Range("U1").Value = "Topic_19"
$ws.Range("A2").Value = "antena3_2019 09 14_morning_new"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.958
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0.042
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("A3").Value = "antena3_2019 09 14_afternoon_new"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("A4").Value = "antena3_2019 09 15_morning_new"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.716
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.284
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("A5").Value = "antena3_2019 09 15_afternoon_new"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.533
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0.467
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("A6").Value = "antena3_2019 09 18_afternoon_new"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 0.736
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0.234
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0.014
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0.015
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0

# --- Update worksheet "Cluster3" ---
$ws = $wb.Worksheets.Item("Cluster3")
$ws.Range("B1").Value = "clusters"
$ws.Range("C1").Value = "Topic_1"
$ws.Range("D1").Value = "Topic_2"
$ws.Range("E1").Value = "Topic_3"
$ws.Range("F1").Value = "Topic_4"
$ws.Range("G1").Value = "Topic_5"
$ws.Range("H1").Value = "Topic_6"
$ws.Range("I1").Value = "Topic_7"
$ws.Range("J1").Value = "Topic_8"
$ws.Range("K1").Value = "Topic_9"
$ws.Range("L1").Value = "Topic_10"
$ws.Range("M1").Value = "Topic_11"
$ws.Range("N1").Value = "Topic_12"
$ws.Range("O1").Value = "Topic_13"
$ws.Range("P1").Value = "Topic_14"
$ws.Range("Q1").Value = "Topic_15"
$ws.Range("R1").Value = "Topic_16"
$ws.Range("S1").Value = "Topic_17"
$ws.Range("T1").Value = "Topic_18"
$ws.Range("U1").Value = "Topic_19"
$ws.Range("A2").Value = "antena3_2019 09 16_afternoon_new"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.207
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0.792
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0

